$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(9)
